$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13c"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.651738666666667
$ws.Range("H2").Value = 7.955216
$ws.Range("I2").Value = 0.5604432756981275
$ws.Range("J2").Value = 0.5604432756981274
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.2138055
$ws.Range("N2").Value = 0.427611
$ws.Range("O2").Value = 0.2116309566486094
$ws.Range("P2").Value = 0.1517955144088434
$ws.Range("Q2").Value = 0.566956311496
$ws.Range("R2").Value = 3.401737868976
$ws.Range("S2").Value = 0.1186071465832751
$ws.Range("T2").Value = 0.08507277533157447

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13c"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.651738666666667
$ws.Range("H3").Value = 7.955216
$ws.Range("I3").Value = 0.5604432756981275
$ws.Range("J3").Value = 0.5604432756981274
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7964696666666667
$ws.Range("N3").Value = 2.389409
$ws.Range("O3").Value = 0.7883690433513906
$ws.Range("P3").Value = 0.8482044855911566
$ws.Range("Q3").Value = 2.112029411927111
$ws.Range("R3").Value = 19.008264707344
$ws.Range("S3").Value = 0.4418361291148524
$ws.Range("T3").Value = 0.4753705003665529

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13c"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5155733333333333
$ws.Range("H4").Value = 1.54672
$ws.Range("I4").Value = 0.10896609512398
$ws.Range("J4").Value = 0.10896609512398
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.2138055
$ws.Range("N4").Value = 0.427611
$ws.Range("O4").Value = 0.2116309566486094
$ws.Range("P4").Value = 0.1517955144088434
$ws.Range("Q4").Value = 0.11023241432
$ws.Range("R4").Value = 0.6613944859199999
$ws.Range("S4").Value = 0.02306059895335126
$ws.Range("T4").Value = 0.0165405644624675

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Tnfsf13b"
$ws.Range("C5").Value = "Tnfrsf13c"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5155733333333333
$ws.Range("H5").Value = 1.54672
$ws.Range("I5").Value = 0.10896609512398
$ws.Range("J5").Value = 0.10896609512398
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7964696666666667
$ws.Range("N5").Value = 2.389409
$ws.Range("O5").Value = 0.7883690433513906
$ws.Range("P5").Value = 0.8482044855911566
$ws.Range("Q5").Value = 0.4106385209422223
$ws.Range("R5").Value = 3.69574668848
$ws.Range("S5").Value = 0.08590549617062875
$ws.Range("T5").Value = 0.09242553066151248

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Tnfsf13b"
$ws.Range("C6").Value = "Tnfrsf13c"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.564190333333333
$ws.Range("H6").Value = 4.692571
$ws.Range("I6").Value = 0.3305906291778926
$ws.Range("J6").Value = 0.3305906291778925
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.2138055
$ws.Range("N6").Value = 0.427611
$ws.Range("O6").Value = 0.2116309566486094
$ws.Range("P6").Value = 0.1517955144088434
$ws.Range("Q6").Value = 0.3344324963135
$ws.Range("R6").Value = 2.006594977881
$ws.Range("S6").Value = 0.0699632111119831
$ws.Range("T6").Value = 0.05018217461480138

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Tnfsf13b"
$ws.Range("C7").Value = "Tnfrsf13c"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.564190333333333
$ws.Range("H7").Value = 4.692571
$ws.Range("I7").Value = 0.3305906291778926
$ws.Range("J7").Value = 0.3305906291778925
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7964696666666667
$ws.Range("N7").Value = 2.389409
$ws.Range("O7").Value = 0.7883690433513906
$ws.Range("P7").Value = 0.8482044855911566
$ws.Range("Q7").Value = 1.245830153393222
$ws.Range("R7").Value = 11.212471380539
$ws.Range("S7").Value = 0.2606274180659094
$ws.Range("T7").Value = 0.2804084545630911

